# Updates the currentAveragePrice*/LevePrice*/LeveProfit* figures (columns H-N) for a
# batch of Leve rows across the ALC, ARM, BSM, GSM and WVR sheets, matching the latest
# scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush / Growth Formula Alpha
$ws.Range("H4").Value = 579.2143
$ws.Range("J4").Value = 397.7143
$ws.Range("L4").Value = 397.7143
$ws.Range("N4").Value = -625.7143
# Row 10: A Jawbreaking Weapon of Staggering Weight / Whispering Maple Wand
$ws.Range("H10").Value = 5313.143
$ws.Range("J10").Value = 5313.143
$ws.Range("L10").Value = 5313.143
$ws.Range("N10").Value = -5899.143
# Row 20: Shut Up and Take My Gil / Ash Wand
$ws.Range("H20").Value = 1850
$ws.Range("I20").Value = 1850
$ws.Range("K20").Value = 1850
$ws.Range("M20").Value = -1620
# Row 21: Book and a Hard Place / Engraved Hard Leather Grimoire
$ws.Range("H21").Value = 3599
$ws.Range("I21").Value = 2999
$ws.Range("J21").Value = 3999
$ws.Range("K21").Value = 2999
$ws.Range("L21").Value = 3999
$ws.Range("M21").Value = -2531
$ws.Range("N21").Value = -4935
# Row 23: There's Something about Bury / Hard Leather Grimoire
$ws.Range("H23").Value = 3599
$ws.Range("I23").Value = 2999
$ws.Range("J23").Value = 3999
$ws.Range("K23").Value = 2999
$ws.Range("L23").Value = 3999
$ws.Range("M23").Value = -2765
$ws.Range("N23").Value = -4467
# Row 35: Conspicuous Conjuration / Whispering Ash Wand
$ws.Range("H35").Value = 1850
$ws.Range("I35").Value = 1850
$ws.Range("K35").Value = 1850
$ws.Range("M35").Value = -1471
# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 1666.1904
$ws.Range("I62").Value = 1611.6666
$ws.Range("K62").Value = 1611.6666
$ws.Range("M62").Value = -987.6666
# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 1666.1904
$ws.Range("I65").Value = 1611.6666
$ws.Range("K65").Value = 8058.333000000001
$ws.Range("M65").Value = -4938.333000000001
# Row 70: Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 1350.3334
$ws.Range("I70").Value = 1200.6666
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 3601.9998
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -3331.9998
$ws.Range("N70").Value = -5040
# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 1350.3334
$ws.Range("I73").Value = 1200.6666
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 3601.9998
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -2665.9998
$ws.Range("N73").Value = -6372
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 4675920.5
$ws.Range("I98").Value = 52808.57
$ws.Range("J98").Value = 37037704
$ws.Range("K98").Value = 52808.57
$ws.Range("L98").Value = 37037704
$ws.Range("M98").Value = -51310.57
$ws.Range("N98").Value = -37040700
# Row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 1367.5
$ws.Range("I111").Value = 1352
$ws.Range("J111").Value = 1383
$ws.Range("K111").Value = 4056
$ws.Range("L111").Value = 4149
$ws.Range("M111").Value = -989
$ws.Range("N111").Value = -10283
# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 4675920.5
$ws.Range("I122").Value = 52808.57
$ws.Range("J122").Value = 37037704
$ws.Range("K122").Value = 158425.71
$ws.Range("L122").Value = 111113112
$ws.Range("M122").Value = -155975.71
$ws.Range("N122").Value = -111118012
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 1448.881
$ws.Range("I132").Value = 1333.3158
$ws.Range("J132").Value = 2546.75
$ws.Range("K132").Value = 3999.9474
$ws.Range("L132").Value = 7640.25
$ws.Range("M132").Value = -1469.9474
$ws.Range("N132").Value = -12700.25
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2353.1538
$ws.Range("I138").Value = 1354.091
$ws.Range("J138").Value = 3646.0588
$ws.Range("K138").Value = 4062.273
$ws.Range("L138").Value = 10938.1764
$ws.Range("M138").Value = 1077.727
$ws.Range("N138").Value = -21218.1764

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 420527.3
$ws.Range("I32").Value = 5057.1294
$ws.Range("J32").Value = 1916219.9
$ws.Range("K32").Value = 5057.1294
$ws.Range("L32").Value = 1916219.9
$ws.Range("M32").Value = -4770.1294
$ws.Range("N32").Value = -1916793.9
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2658.6956
$ws.Range("I61").Value = 2926.9412
$ws.Range("J61").Value = 1898.6666
$ws.Range("K61").Value = 2926.9412
$ws.Range("L61").Value = 1898.6666
$ws.Range("M61").Value = -2714.9412
$ws.Range("N61").Value = -2322.6666
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 4059.4856
$ws.Range("I122").Value = 4077.9033
$ws.Range("K122").Value = 12233.7099
$ws.Range("M122").Value = -9783.7099
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2658.6956
$ws.Range("I136").Value = 2926.9412
$ws.Range("J136").Value = 1898.6666
$ws.Range("K136").Value = 8780.8236
$ws.Range("L136").Value = 5695.9998
$ws.Range("M136").Value = -6230.8236
$ws.Range("N136").Value = -10795.9998

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 1165.75
$ws.Range("I94").Value = 872.4375
$ws.Range("J94").Value = 2339
$ws.Range("K94").Value = 872.4375
$ws.Range("L94").Value = 2339
$ws.Range("M94").Value = -421.4375
$ws.Range("N94").Value = -3241

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 1852.6818
$ws.Range("I97").Value = 1201.5385
$ws.Range("J97").Value = 2793.2222
$ws.Range("K97").Value = 1201.5385
$ws.Range("L97").Value = 2793.2222
$ws.Range("M97").Value = -705.5385000000001
$ws.Range("N97").Value = -3785.2222
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1840.0834
$ws.Range("I102").Value = 1887.3529
$ws.Range("J102").Value = 1725.2858
$ws.Range("K102").Value = 1887.3529
$ws.Range("L102").Value = 1725.2858
$ws.Range("M102").Value = -265.3529000000001
$ws.Range("N102").Value = -4969.2858
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 1916.8541
$ws.Range("I132").Value = 1537.4634
$ws.Range("J132").Value = 4139
$ws.Range("K132").Value = 4612.3902
$ws.Range("L132").Value = 12417
$ws.Range("M132").Value = -2082.3902
$ws.Range("N132").Value = -17477

$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables / Hempen Underpants
$ws.Range("H2").Value = 650.3333
$ws.Range("I2").Value = 475.5
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 475.5
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -363.5
$ws.Range("N2").Value = -1224
# Row 14: Hat in Hand / Straw Hat
$ws.Range("H14").Value = 111112650
$ws.Range("J14").Value = 2640
$ws.Range("L14").Value = 2640
$ws.Range("N14").Value = -2976
# Row 15: Workplace Safety / Cotton Scarf
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
# Row 19: Dirt Cheap / Stablehand's Hat
$ws.Range("H19").Value = 1975
$ws.Range("J19").Value = 1975
$ws.Range("L19").Value = 1975
$ws.Range("N19").Value = -2323
# Row 29: Getting Handsy / Cotton Dress Gloves
$ws.Range("H29").Value = 6667100
$ws.Range("I29").Value = 6667100
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 6667100
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -6666810
$ws.Range("N29").ClearContents()
